$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "active_power_im" in C1 and value 987654321 in C2
$ws.Range("C1").Value = "active_power_im"
$ws.Range("C2").Value = 987654321

# Update the active selection to C7 (mirrors the diff's selection change)
$ws.Range("C7").Select()
